$d = $word.ActiveDocument

# Delete the large block of paragraphs that make up the entire
# "thomas-satSef" section: the page-break paragraph through the final
# "utils" paragraph at the end of the document.
$startPara = $d.Paragraphs.Item(22)
$endPara = $d.Paragraphs.Item(32)
$r = $d.Range($startPara.Range.Start, $endPara.Range.End)
$r.Delete()

# Delete the "setSatPaper" and "vss2020Fig" paragraphs (including their
# paragraph marks) from the remaining "paperFigs" listing.
$vss2020Fig = $d.Paragraphs.Item(16)
$vss2020Fig.Range.Delete()

$setSatPaper = $d.Paragraphs.Item(15)
$setSatPaper.Range.Delete()
